$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 9.1 = 37060.17 pesos", "1000 Bs = 9.18 = 37466.48 pesos")
$text = $text.Replace("37060.17 pesos = 9.07 = 951.93 Bs", "37466.48 pesos = 9.12 = 927.09 Bs")
$cell.Value = $text

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 108.9
$wsTasas.Range("O10").Value = 4080.1
$wsTasas.Range("N12").Value = 4110
$wsTasas.Range("O12").Value = 101.7
